{"js": "// Fix a typo in the blog post paragraph: \"Do Not Repeat (yourself)\" -> \"Don't Repeat Yourself\"\nconst oldPhrase = \"Do Not Repeat (yourself)\";\nconst newPhrase = \"Don't Repeat Yourself\";\n\nconst results = context.document.body.search(oldPhrase, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found: \" + oldPhrase);\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newPhrase, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix a typo in the blog post paragraph: \"Do Not Repeat (yourself)\" -> \"Don't Repeat Yourself\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Do Not Repeat (yourself)\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1   # wdFindContinue\n\n$found = $find.Execute()\n\nif ($found) {\n    # Assign directly to the matched Range's Text (instead of Find's\n    # Replacement mechanism) so the straight apostrophe in the replacement\n    # is preserved verbatim instead of being smart-quoted.\n    $rng = $find.Parent\n    $rng.Text = \"Don't Repeat Yourself\"\n} else {\n    throw \"Target phrase not found: Do Not Repeat (yourself)\"\n}\n"}
